$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row cells: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range into an Excel Table ("Table1")
$usedRange = $ws.Range("A1:U71")
$listObject = $ws.ListObjects.Add(1, $usedRange, 0, 1)
$listObject.Name = "Table1"

# Re-assert the column header names (ListObjects.Add may otherwise reuse
# whatever text was already present) so the table column names match the
# cell values set above.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $listObject.ListColumns.Item($i + 1).Name = $headers[$i]
}

# 3) Freeze the header row (split after row 1, keep viewing from A2)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
